$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (GitHub Actions data pull),
# including the Polkadot / WrappedEther rows trading places (rows 12-13).
#
# Every target cell here is plain Text in the source sheet (inlineStr), and several
# of the new values look like numbers/dates (e.g. "20.40", "1.003"). Assigning those
# bare via .Value lets Excel auto-coerce them to Number (dropping the trailing zero,
# i.e. 20.40 -> 20.4) which would not match the source formatting. Prefixing with an
# apostrophe forces Text, exactly like the original cells, without the quote itself
# becoming part of the stored value.

# Row 2
$ws.Range("D2").Value = '''25.891.45'
$ws.Range("E2").Value = '''  +0.23%  '

# Row 3
$ws.Range("D3").Value = '''1.638.06'
$ws.Range("E3").Value = '''  +0.51%  '

# Row 4
$ws.Range("E4").Value = '''  +0.20%  '

# Row 5
$ws.Range("D5").Value = '''215.42'
$ws.Range("E5").Value = '''  +0.13%  '

# Row 6
$ws.Range("D6").Value = '''0.5093'
$ws.Range("E6").Value = '''  -0.40%  '

# Row 7
$ws.Range("D7").Value = '''1.003'
$ws.Range("E7").Value = '''  +0.26%  '

# Row 8
$ws.Range("D8").Value = '''0.2590'
$ws.Range("E8").Value = '''  +0.81%  '

# Row 9
$ws.Range("D9").Value = '''0.06439'
$ws.Range("E9").Value = '''  +1.50%  '

# Row 10
$ws.Range("D10").Value = '''20.40'
$ws.Range("E10").Value = '''  +4.65%  '

# Row 11
$ws.Range("D11").Value = '''0.07819'
$ws.Range("E11").Value = '''  +0.44%  '

# Row 12
$ws.Range("B12").Value = '''WrappedEther'
$ws.Range("C12").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '''1.661.18'
$ws.Range("E12").Value = '''  +1.76%  '

# Row 13
$ws.Range("B13").Value = '''Polkadot'
$ws.Range("C13").Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.278'
$ws.Range("E13").Value = '''  +0.56%  '

# Row 14
$ws.Range("D14").Value = '''1.865.62'
$ws.Range("E14").Value = '''  +0.74%  '

# Row 15
$ws.Range("D15").Value = '''0.5606'
$ws.Range("E15").Value = '''  +1.35%  '

# Row 16
$ws.Range("D16").Value = '''0.0₅7682'
$ws.Range("E16").Value = '''  +2.26%  '

# Row 17
$ws.Range("D17").Value = '''63.35'
$ws.Range("E17").Value = '''  -0.48%  '

# Row 18
$ws.Range("D18").Value = '''25.901.43'
$ws.Range("E18").Value = '''  +0.27%  '

# Row 19
$ws.Range("D19").Value = '''1.003'
$ws.Range("E19").Value = '''  +0.22%  '

# Row 20
$ws.Range("D20").Value = '''194.36'
$ws.Range("E20").Value = '''  -0.21%  '

# Row 21
$ws.Range("D21").Value = '''4.391'
$ws.Range("E21").Value = '''  -0.95%  '

# Row 22
$ws.Range("D22").Value = '''9.959'
$ws.Range("E22").Value = '''  +1.67%  '

# Row 23
$ws.Range("D23").Value = '''6.157'
$ws.Range("E23").Value = '''  +2.29%  '

# Row 24
$ws.Range("D24").Value = '''1.003'
$ws.Range("E24").Value = '''  +0.23%  '

# Row 25
$ws.Range("D25").Value = '''1.791'
$ws.Range("E25").Value = '''  -5.05%  '

# Row 26
$ws.Range("D26").Value = '''138.06'
$ws.Range("E26").Value = '''  -2.43%  '

# Row 27
$ws.Range("E27").Value = '''  -1.53%  '

# Row 28
$ws.Range("D28").Value = '''6.844'
$ws.Range("E28").Value = '''  +1.80%  '

# Row 29
$ws.Range("E29").Value = '''  -0.03%  '

# Row 30
$ws.Range("D30").Value = '''1.243'
$ws.Range("E30").Value = '''  +0.33%  '

# Row 31
$ws.Range("D31").Value = '''0.04976'
$ws.Range("E31").Value = '''  +2.16%  '

# Row 32
$ws.Range("D32").Value = '''3.304'
$ws.Range("E32").Value = '''  +1.31%  '

# Row 33
$ws.Range("D33").Value = '''3.253'
$ws.Range("E33").Value = '''  +2.56%  '

# Row 34
$ws.Range("E34").Value = '''  +1.77%  '

# Row 35
$ws.Range("D35").Value = '''2.386'
$ws.Range("E35").Value = '''  +1.14%  '

# Row 36
$ws.Range("D36").Value = '''0.9048'
$ws.Range("E36").Value = '''  +0.91%  '

# Row 37
$ws.Range("E37").Value = '''  +1.44%  '

# Row 38
$ws.Range("E38").Value = '''  +0.75%  '

# Row 39
$ws.Range("D39").Value = '''1.138.07'
$ws.Range("E39").Value = '''  +1.89%  '

# Row 40
$ws.Range("E40").Value = '''  +1.61%  '

# Row 41
$ws.Range("D41").Value = '''1.003'
$ws.Range("E41").Value = '''  +0.44%  '

# Row 42
$ws.Range("D42").Value = '''99.60'
$ws.Range("E42").Value = '''  +2.16%  '

# Row 43
$ws.Range("D43").Value = '''5.477'
$ws.Range("E43").Value = '''  -1.31%  '

# Row 44
$ws.Range("D44").Value = '''0.8023'
$ws.Range("E44").Value = '''  +0.56%  '

# Row 45
$ws.Range("D45").Value = '''0.0₈114'
$ws.Range("E45").Value = '''  -1.44%  '

# Row 46
$ws.Range("E46").Value = '''  +1.77%  '

# Row 47
$ws.Range("D47").Value = '''0.4253'
$ws.Range("E47").Value = '''  -3.88%  '

# Row 48
$ws.Range("D48").Value = '''7.779'
$ws.Range("E48").Value = '''  +2.08%  '

# Row 49
$ws.Range("D49").Value = '''0.05065'
$ws.Range("E49").Value = '''  -1.11%  '

# Row 50
$ws.Range("E50").Value = '''  +0.03%  '

# Row 51
$ws.Range("E51").Value = '''  +0.27%  '
